# Fix conditional forecasts so wages feed into shortages
# Update the "shortage" worksheet summary stats with corrected values.
# (Scientific notation literals are not supported by the script parser,
# so all values are written out in plain decimal form.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shortage")

# Row 2: l1.shortage through l4.shortage
$ws.Range("B2").Value = 0.5934328645957652
$ws.Range("C2").Value = 0.00001364928244292737
$ws.Range("D2").Value = 0.0001104115643712736

# Row 3: excess_demand through l4.excess_demand
$ws.Range("B3").Value = 15.376022482534220
$ws.Range("C3").Value = 0.1219970040973848
$ws.Range("D3").Value = 0.00002330505244000491

# Row 4: gscpi through l4.gscpi
$ws.Range("B4").Value = 3.451610643570803
$ws.Range("C4").Value = 0.0001177553333192535
$ws.Range("D4").Value = 0.000001757868254029659

# Row 6: Long-run excess demand multiplier
$ws.Range("B6").Value = 37.819147549263668

# Row 7: Long-run GSCPI multiplier
$ws.Range("B7").Value = 8.489644988493703
